$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 195.36363
$ws.Range("I9").Value = 203.875
$ws.Range("K9").Value = 203.875
$ws.Range("M9").Value = -34.875
$ws.Range("H28").Value = 555.1818
$ws.Range("J28").Value = 1199.5
$ws.Range("L28").Value = 1199.5
$ws.Range("N28").Value = -2169.5
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -2719
$ws.Range("H87").Value = 56608
$ws.Range("J87").Value = 56608
$ws.Range("L87").Value = 56608
$ws.Range("N87").Value = -59104
$ws.Range("H90").Value = 56608
$ws.Range("J90").Value = 56608
$ws.Range("L90").Value = 169824
$ws.Range("N90").Value = -182304
$ws.Range("H135").Value = 2535.9092
$ws.Range("I135").Value = 1766.7778
$ws.Range("K135").Value = 15901.0002
$ws.Range("M135").Value = -13366.0002
$ws.Range("H137").Value = 9820.474
$ws.Range("I137").Value = 1626.909
$ws.Range("K137").Value = 4880.727000000001
$ws.Range("M137").Value = -2330.727000000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1900
$ws.Range("I2").Value = 1900
$ws.Range("K2").Value = 1900
$ws.Range("M2").Value = -1787
$ws.Range("H45").Value = 2550.9167
$ws.Range("I45").Value = 1601.375
$ws.Range("K45").Value = 1601.375
$ws.Range("M45").Value = -1224.375
$ws.Range("H102").Value = 2784.8333
$ws.Range("I102").Value = 2965.2727
$ws.Range("K102").Value = 2965.2727
$ws.Range("M102").Value = -1343.2727
$ws.Range("H116").Value = 1900
$ws.Range("I116").Value = 1900
$ws.Range("K116").Value = 1900
$ws.Range("M116").Value = 394
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H132").Value = 522050.9
$ws.Range("I132").Value = 556578.9
$ws.Range("J132").Value = 4131.6665
$ws.Range("K132").Value = 1669736.7
$ws.Range("L132").Value = 12394.9995
$ws.Range("M132").Value = -1667206.7
$ws.Range("N132").Value = -17454.9995

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1900
$ws.Range("I3").Value = 1900
$ws.Range("K3").Value = 1900
$ws.Range("M3").Value = -1786
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0
$ws.Range("H101").Value = 62500
$ws.Range("J101").Value = 62500
$ws.Range("L101").Value = 62500
$ws.Range("N101").Value = -68990
$ws.Range("H103").Value = 38665
$ws.Range("J103").Value = 38665
$ws.Range("L103").Value = 38665
$ws.Range("N103").Value = -41009
$ws.Range("H134").Value = 5933.968
$ws.Range("I134").Value = 2793.1428
$ws.Range("K134").Value = 8379.428400000001
$ws.Range("M134").Value = -5844.428400000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2462.85
$ws.Range("I31").Value = 2403.2222
$ws.Range("J31").Value = 2999.5
$ws.Range("K31").Value = 2403.2222
$ws.Range("L31").Value = 2999.5
$ws.Range("M31").Value = -2108.2222
$ws.Range("N31").Value = -3589.5
$ws.Range("H34").Value = 2462.85
$ws.Range("I34").Value = 2403.2222
$ws.Range("J34").Value = 2999.5
$ws.Range("K34").Value = 2403.2222
$ws.Range("L34").Value = 2999.5
$ws.Range("M34").Value = -2201.2222
$ws.Range("N34").Value = -3403.5
$ws.Range("H62").Value = 3333.3333
$ws.Range("I62").Value = 3250
$ws.Range("K62").Value = 3250
$ws.Range("M62").Value = -2626
$ws.Range("H65").Value = 3333.3333
$ws.Range("I65").Value = 3250
$ws.Range("K65").Value = 16250
$ws.Range("M65").Value = -13130
$ws.Range("H122").Value = 16661.334
$ws.Range("I122").Value = 3018.1667
$ws.Range("K122").Value = 9054.500100000001
$ws.Range("M122").Value = -6604.500100000001
$ws.Range("H132").Value = 2294.6843
$ws.Range("I132").Value = 2240.7646
$ws.Range("J132").Value = 2753
$ws.Range("K132").Value = 6722.293799999999
$ws.Range("L132").Value = 8259
$ws.Range("M132").Value = -4192.293799999999
$ws.Range("N132").Value = -13319

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2664.4707
$ws.Range("I3").Value = 2468.5
$ws.Range("J3").Value = 5800
$ws.Range("K3").Value = 7405.5
$ws.Range("L3").Value = 17400
$ws.Range("M3").Value = -7293.5
$ws.Range("N3").Value = -17624
$ws.Range("H75").Value = 5015
$ws.Range("J75").Value = 5015
$ws.Range("L75").Value = 15045
$ws.Range("N75").Value = -17041
$ws.Range("H78").Value = 5015
$ws.Range("J78").Value = 5015
$ws.Range("L78").Value = 45135
$ws.Range("N78").Value = -55119
$ws.Range("H114").Value = 4210.923
$ws.Range("I114").Value = 2083.3333
$ws.Range("J114").Value = 4849.2
$ws.Range("K114").Value = 6249.999899999999
$ws.Range("L114").Value = 14547.6
$ws.Range("M114").Value = -2995.999899999999
$ws.Range("N114").Value = -21055.6
$ws.Range("H122").Value = 1613831.8
$ws.Range("I122").Value = 16129032
$ws.Range("J122").Value = 1031.6666
$ws.Range("K122").Value = 145161288
$ws.Range("L122").Value = 9284.999400000001
$ws.Range("M122").Value = -145158838
$ws.Range("N122").Value = -14184.9994
$ws.Range("H131").Value = 1728.5625
$ws.Range("I131").Value = 904.75
$ws.Range("J131").Value = 1846.25
$ws.Range("K131").Value = 2714.25
$ws.Range("L131").Value = 5538.75
$ws.Range("M131").Value = 2325.75
$ws.Range("N131").Value = -15618.75
$ws.Range("H132").Value = 523.6
$ws.Range("I132").Value = 576.5714
$ws.Range("K132").Value = 5189.1426
$ws.Range("M132").Value = -2659.1426
$ws.Range("H133").Value = 7792.909
$ws.Range("I133").Value = 7872.2
$ws.Range("K133").Value = 23616.6
$ws.Range("M133").Value = -18556.6

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 17299.666
$ws.Range("J20").Value = 10949.5
$ws.Range("L20").Value = 10949.5
$ws.Range("N20").Value = -11439.5
$ws.Range("H24").Value = 10666.667
$ws.Range("J24").Value = 10666.667
$ws.Range("L24").Value = 10666.667
$ws.Range("N24").Value = -11012.667
$ws.Range("H122").Value = 75078.28999999999
$ws.Range("I122").Value = 94236.09
$ws.Range("K122").Value = 282708.27
$ws.Range("M122").Value = -280258.27
$ws.Range("H132").Value = 13680.566
$ws.Range("I132").Value = 14789.704
$ws.Range("K132").Value = 44369.112
$ws.Range("M132").Value = -41839.112

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1700
$ws.Range("I68").Value = 1700
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1700
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -951
$ws.Range("H71").Value = 1700
$ws.Range("I71").Value = 1700
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 8500
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -4756
$ws.Range("H122").Value = 8498.75
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 8498.75
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 25496.25
$ws.Range("N122").Value = -30396.25
$ws.Range("H132").Value = 3574170.5
$ws.Range("I132").Value = 4548099
$ws.Range("K132").Value = 13644297
$ws.Range("M132").Value = -13641767
$ws.Range("H136").Value = 12682.546
$ws.Range("I136").Value = 4519.5
$ws.Range("K136").Value = 13558.5
$ws.Range("M136").Value = -11008.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 20333.334
$ws.Range("I31").Value = 8000
$ws.Range("J31").Value = 45000
$ws.Range("K31").Value = 8000
$ws.Range("L31").Value = 45000
$ws.Range("M31").Value = -7652
$ws.Range("N31").Value = -45696
$ws.Range("H46").Value = 60995
$ws.Range("J46").Value = 60995
$ws.Range("L46").Value = 60995
$ws.Range("N46").Value = -61457
$ws.Range("H81").Value = 3009.7144
$ws.Range("J81").Value = 4709.25
$ws.Range("L81").Value = 9418.5
$ws.Range("N81").Value = -11540.5
$ws.Range("H84").Value = 3009.7144
$ws.Range("J84").Value = 4709.25
$ws.Range("L84").Value = 47092.5
$ws.Range("N84").Value = -57700.5
$ws.Range("H122").Value = 61530.844
$ws.Range("I122").Value = 2464.4443
$ws.Range("K122").Value = 7393.3329
$ws.Range("M122").Value = -4943.3329
$ws.Range("H126").Value = 2233.9
$ws.Range("I126").Value = 2233.9
$ws.Range("K126").Value = 6701.700000000001
$ws.Range("M126").Value = -4231.700000000001
$ws.Range("H132").Value = 5090.75
$ws.Range("I132").Value = 5090.75
$ws.Range("K132").Value = 15272.25
$ws.Range("M132").Value = -12742.25
$ws.Range("H134").Value = 60995
$ws.Range("J134").Value = 60995
$ws.Range("L134").Value = 182985
$ws.Range("N134").Value = -188055
$ws.Range("H136").Value = 1250.5385
$ws.Range("I136").Value = 1250.7273
$ws.Range("K136").Value = 3752.1819
$ws.Range("M136").Value = -1202.1819
